# Append "Request" / "Response" to the Common Name (column A) display names,
# based on whether the corresponding Document Type Identifier (column B)
# refers to a Request or a Response.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Text
    $urn = $ws.Cells.Item($r, 2).Text

    if ([string]::IsNullOrEmpty($name) -or [string]::IsNullOrEmpty($urn)) {
        continue
    }

    if ($urn -like "*::Request##*") {
        $ws.Cells.Item($r, 1).Value = "$name Request"
    }
    elseif ($urn -like "*::Response##*") {
        $ws.Cells.Item($r, 1).Value = "$name Response"
    }
}

# Move the active selection to A20, matching the post-edit cursor position.
$ws.Range("A20").Select() | Out-Null
